$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 278
$ws.Range("F4").Value = 286
$ws.Range("F5").Value = 2905
$ws.Range("F8").Value = 2271
$ws.Range("F9").Value = 1525
$ws.Range("F10").Value = 44
$ws.Range("F12").Value = 93
$ws.Range("F13").Value = 2601
$ws.Range("F15").Value = 1448
$ws.Range("F16").Value = 6445
$ws.Range("F18").Value = 6214
$ws.Range("F20").Value = 2158
$ws.Range("F21").Value = 2991
$ws.Range("F22").Value = 3414
$ws.Range("F23").Value = 200
$ws.Range("F24").Value = 22
$ws.Range("F25").Value = 1692
$ws.Range("F26").Value = 60
$ws.Range("F27").Value = 281
$ws.Range("F29").Value = 156
$ws.Range("F30").Value = 17
$ws.Range("F31").Value = 350
$ws.Range("F32").Value = 1066
$ws.Range("F33").Value = 2305
$ws.Range("F34").Value = 8
$ws.Range("F35").Value = 145
$ws.Range("F36").Value = 326
$ws.Range("F37").Value = 878
$ws.Range("F38").Value = 178
$ws.Range("F39").Value = 415
$ws.Range("F40").Value = 478

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 108
$ws.Range("F6").Value = 64
$ws.Range("F10").Value = 94
$ws.Range("F13").Value = 201
$ws.Range("F20").Value = 7
$ws.Range("F21").Value = 53

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 108
$ws.Range("F4").Value = 64
$ws.Range("F5").Value = 278
$ws.Range("F7").Value = 286
$ws.Range("F9").Value = 2905
$ws.Range("F11").Value = 2271
$ws.Range("F12").Value = 1525
$ws.Range("F13").Value = 44
$ws.Range("F15").Value = 93
$ws.Range("F17").Value = 2601
$ws.Range("F18").Value = 1448
$ws.Range("F19").Value = 201
$ws.Range("F23").Value = 6445
$ws.Range("F25").Value = 6214
$ws.Range("F26").Value = 2158
$ws.Range("F27").Value = 2991
$ws.Range("F28").Value = 3414
$ws.Range("F30").Value = 200
$ws.Range("F33").Value = 1692
$ws.Range("F34").Value = 7
$ws.Range("F35").Value = 53
$ws.Range("F36").Value = 281
$ws.Range("F38").Value = 156
$ws.Range("F39").Value = 17
$ws.Range("F40").Value = 350
$ws.Range("F42").Value = 2305
$ws.Range("F43").Value = 8
$ws.Range("F44").Value = 145
$ws.Range("F45").Value = 326
$ws.Range("F46").Value = 878
$ws.Range("F47").Value = 178
$ws.Range("F48").Value = 415
$ws.Range("F49").Value = 478

